$wb = $excel.ActiveWorkbook

# Update "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1449
$ws1.Range("F4").Value = 96
$ws1.Range("F6").Value = 19

# Update "全部类型" sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1449
$ws4.Range("F4").Value = 96
$ws4.Range("F6").Value = 19
